$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.112518548965454
$ws.Range("B1").Value = 2.482971668243408
$ws.Range("C1").Value = 2.493942975997925
$ws.Range("D1").Value = 2.860424518585205
$ws.Range("E1").Value = 0.5695806741714478
